$p = $ppt.ActivePresentation

# 1. Slide 10 ("Exercise 1" divider slide): append a new paragraph.
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(1)
$sh10.TextFrame.TextRange.Text = "Exercise 1`rHands On Exercise 1 Part A"

# 2. Slide 14 ("Exercise 2" divider slide): duplicate it first (so the new
#    slide inherits the original single-paragraph "Exercise 2" shape/creationId),
#    landing right after it as the new slide 15, then append the new paragraph
#    to the original slide 14.
$s14 = $p.Slides.Item(14)
$newSlide = $s14.Duplicate()

$sh14 = $s14.Shapes.Item(1)
$sh14.TextFrame.TextRange.Text = "Exercise 2`rHands On Exercise 1 Part B"

# 3. The duplicated slide becomes the new slide 15 ("Exercise 3 (optional)").
$shNew = $newSlide.Shapes.Item(1)
$shNew.TextFrame.TextRange.Text = "Exercise 3 (optional)`rHands On Exercise 1 Part C"
